$d = $word.ActiveDocument

# Locate the three target paragraphs (the "Undefined variables..." entry plus the
# two trailing highlighted TODO entries) by their distinctive text.
$targetCount = $d.Paragraphs.Count
$firstIdx = -1
$lastIdx = -1
for ($i = 1; $i -le $targetCount; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Undefined variables*") {
        $firstIdx = $i
    }
    if ($t -like "In Profile p0age*") {
        $lastIdx = $i
    }
}

$startPos = $d.Paragraphs.Item($firstIdx).Range.Start
$endPos = $d.Paragraphs.Item($lastIdx).Range.End

# Remove the old "Undefined variables ..." paragraph plus the two stray
# highlighted TODO paragraphs that trail the list (they get reworked/replaced
# below, with the bookmark relocated to a new trailing empty paragraph).
$oldRange = $d.Range($startPos, $endPos)
$oldRange.Delete()

# Re-insert the rewritten paragraph plus all of the new TODO-list entries as
# raw OOXML so the run/paragraph split matches exactly (highlight removed,
# times updated, new items appended, bookmark moved to a fresh final
# paragraph).
$insertXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rStyle w:val="Strong"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Undefined </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>variables</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> images in edit tenant page-</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> 18.05 </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>–</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>11.34</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> - done</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rStyle w:val="Strong"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>Give full name in new bill page in name columns drop down list – 13.21 – 13.42 (pending)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rStyle w:val="Strong"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">In Profile page – put image size same for all images – </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>done</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rStyle w:val="Strong"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>Make landlord expenses lists</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> - done</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rStyle w:val="Strong"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>confirmation before delete</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>on all page</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">– 14.24 </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>–</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>15.</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>44</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rStyle w:val="Strong"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>show active page in sidebar options</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>–</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>15.44</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>–</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>16.19</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rStyle w:val="Strong"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>in edit pages show last value even in drop down list menus</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> – 16.19 </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>–</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rStyle w:val="Strong"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">when I click on back button success message appear in </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>laravel</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rStyle w:val="Strong"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>user image show</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>in different palces</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">- </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>done</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rStyle w:val="Strong"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>give fix size for users images</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> - done</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rStyle w:val="Strong"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">in edit old value show of image </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>–</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> pending</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rStyle w:val="Strong"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">active page show when click on menu </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>–</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> pending</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rStyle w:val="Strong"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">update image in tenant list </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>–</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> pending</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rStyle w:val="Strong"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>null room show in tenant list</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rStyle w:val="Strong"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>I</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">n dashboard every boxes shows list according to their name </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rStyle w:val="Strong"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Make notifications functions in profile </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>dropdown</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> menus</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rStyle w:val="Strong"/><w:sz w:val="24"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$insertPoint = $d.Range($startPos, $startPos)
$insertPoint.InsertXML($insertXml)

# InsertXML does not round-trip the "Strong" character style reference on
# individual runs, so re-apply it explicitly to every run we just inserted
# (this paints each run's Range -- excluding the trailing paragraph mark --
# with the Strong style, which serializes back out as <w:rStyle w:val="Strong"/>).
$newParaCount = $d.Paragraphs.Count - $targetCount + ($lastIdx - $firstIdx + 1)
$firstNewIdx = $firstIdx
$lastNewIdx = $firstIdx + $newParaCount - 1
for ($i = $firstNewIdx; $i -le $lastNewIdx; $i++) {
    $p = $d.Paragraphs.Item($i)
    $r = $p.Range
    if (($r.End - $r.Start) -gt 0) {
        $r.End = $r.End - 1
        if (($r.End - $r.Start) -gt 0) {
            $r.Style = "Strong"
        }
    }
}

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
